$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login email (F5) and integration PID (H5) values
$ws.Range("F5").Value = "yendeli98@gmail.com"
$ws.Range("H5").Value = "4327757"

# Scroll the sheet view so that column D is the left-most visible column
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
